$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pairs of rows whose B:AC (all columns except the leading index column A)
# content must be swapped with each other.
$pairs = @(
    @(644, 645),
    @(694, 695),
    @(701, 702),
    @(716, 717),
    @(719, 720),
    @(740, 741),
    @(777, 778),
    @(783, 784),
    @(786, 787),
    @(799, 800),
    @(811, 812),
    @(842, 843),
    @(894, 895),
    @(1094, 1095),
    @(1120, 1121)
)

foreach ($pair in $pairs) {
    $rowA = $pair[0]
    $rowB = $pair[1]

    $rangeA = $ws.Range("B$rowA`:AC$rowA")
    $rangeB = $ws.Range("B$rowB`:AC$rowB")

    $valA = $rangeA.Value2
    $valB = $rangeB.Value2

    $rangeA.Value2 = $valB
    $rangeB.Value2 = $valA
}

Write-Output "swapped $($pairs.Count) row pairs"
